$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.381.98"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 4).Value = "1.846.03"
$ws.Cells.Item(3, 5).Value = "  -0.14%  "
$ws.Cells.Item(4, 5).Value = "  -0.14%  "
$ws.Cells.Item(5, 4).Value = "'240.88"
$ws.Cells.Item(5, 5).Value = "  -0.94%  "
$ws.Cells.Item(6, 4).Value = "'0.6271"
$ws.Cells.Item(6, 5).Value = "  -3.50%  "
$ws.Cells.Item(7, 5).Value = "  -0.08%  "
$ws.Cells.Item(8, 4).Value = "'0.07592"
$ws.Cells.Item(8, 5).Value = "  +1.28%  "
$ws.Cells.Item(9, 4).Value = "'0.2966"
$ws.Cells.Item(9, 5).Value = "  -0.46%  "
$ws.Cells.Item(10, 4).Value = "'24.41"
$ws.Cells.Item(10, 5).Value = "  -0.07%  "
$ws.Cells.Item(11, 4).Value = "2.111.71"
$ws.Cells.Item(11, 5).Value = "  +13.81%  "
$ws.Cells.Item(12, 4).Value = "'0.07714"
$ws.Cells.Item(12, 5).Value = "  +1.17%  "
$ws.Cells.Item(13, 4).Value = "'0.6870"
$ws.Cells.Item(13, 5).Value = "  +0.39%  "
$ws.Cells.Item(14, 4).Value = "'4.981"
$ws.Cells.Item(14, 5).Value = "  -0.82%  "
$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "2.290.98"
$ws.Cells.Item(15, 5).Value = "  +8.68%  "
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(16, 4).Value = "'82.89"
$ws.Cells.Item(16, 5).Value = "  -0.77%  "
$ws.Cells.Item(17, 4).Value = "'0.000009902"
$ws.Cells.Item(17, 5).Value = "  +4.98%  "
$ws.Cells.Item(18, 4).Value = "'6.151"
$ws.Cells.Item(18, 5).Value = "  +0.73%  "
$ws.Cells.Item(19, 4).Value = "29.417.18"
$ws.Cells.Item(19, 5).Value = "  -0.45%  "
$ws.Cells.Item(20, 4).Value = "'231.29"
$ws.Cells.Item(20, 5).Value = "  -2.90%  "
$ws.Cells.Item(21, 4).Value = "'12.50"
$ws.Cells.Item(21, 5).Value = "  -0.72%  "
$ws.Cells.Item(22, 4).Value = "'0.9999"
$ws.Cells.Item(22, 5).Value = "  -0.06%  "
$ws.Cells.Item(23, 4).Value = "'7.596"
$ws.Cells.Item(23, 5).Value = "  -1.37%  "
$ws.Cells.Item(24, 4).Value = "'1.001"
$ws.Cells.Item(25, 4).Value = "'154.63"
$ws.Cells.Item(25, 5).Value = "  -1.79%  "
$ws.Cells.Item(26, 4).Value = "'0.1387"
$ws.Cells.Item(26, 5).Value = "  -2.10%  "
$ws.Cells.Item(27, 4).Value = "'8.460"
$ws.Cells.Item(27, 5).Value = "  -0.52%  "
$ws.Cells.Item(28, 5).Value = "  -0.87%  "
$ws.Cells.Item(29, 4).Value = "'1.474"
$ws.Cells.Item(29, 5).Value = "  -0.96%  "
$ws.Cells.Item(30, 5).Value = "  -4.40%  "
$ws.Cells.Item(31, 4).Value = "'1.251"
$ws.Cells.Item(31, 5).Value = "  -0.95%  "
$ws.Cells.Item(32, 4).Value = "'4.118"
$ws.Cells.Item(32, 5).Value = "  -0.50%  "
$ws.Cells.Item(33, 4).Value = "'4.016"
$ws.Cells.Item(33, 5).Value = "  -1.24%  "
$ws.Cells.Item(34, 4).Value = "'1.865"
$ws.Cells.Item(34, 5).Value = "  +0.44%  "
$ws.Cells.Item(35, 5).Value = "  -2.18%  "
$ws.Cells.Item(36, 4).Value = "'0.7170"
$ws.Cells.Item(36, 5).Value = "  -1.10%  "
$ws.Cells.Item(37, 4).Value = "'2.595"
$ws.Cells.Item(37, 5).Value = "  +0.05%  "
$ws.Cells.Item(38, 4).Value = "1.247.38"
$ws.Cells.Item(38, 5).Value = "  +3.92%  "
$ws.Cells.Item(39, 4).Value = "'2.794"
$ws.Cells.Item(39, 5).Value = "  -0.05%  "
$ws.Cells.Item(40, 4).Value = "'0.01801"
$ws.Cells.Item(40, 5).Value = "  +1.12%  "
$ws.Cells.Item(41, 4).Value = "2.211.05"
$ws.Cells.Item(41, 5).Value = "  +9.68%  "
$ws.Cells.Item(42, 4).Value = "'0.9076"
$ws.Cells.Item(42, 5).Value = "  -0.04%  "
$ws.Cells.Item(43, 4).Value = "'6.088"
$ws.Cells.Item(43, 5).Value = "  -2.20%  "
$ws.Cells.Item(44, 5).Value = "  -0.11%  "
$ws.Cells.Item(45, 2).Value = "Aave"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(45, 4).Value = "'67.40"
$ws.Cells.Item(45, 5).Value = "  +1.35%  "
$ws.Cells.Item(46, 2).Value = "Quant"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(46, 4).Value = "'101.62"
$ws.Cells.Item(46, 5).Value = "  +0.00%  "
$ws.Cells.Item(47, 4).Value = "'7.302"
$ws.Cells.Item(47, 5).Value = "  -1.57%  "
$ws.Cells.Item(48, 4).Value = "'9.167"
$ws.Cells.Item(48, 5).Value = "  +0.18%  "
$ws.Cells.Item(49, 5).Value = "  -5.44%  "
$ws.Cells.Item(50, 4).Value = "'0.4013"
$ws.Cells.Item(50, 5).Value = "  -0.83%  "
$ws.Cells.Item(51, 4).Value = "'1.695"
$ws.Cells.Item(51, 5).Value = "  +2.82%  "
